$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "321.11"
Set-TextValue "E2" "6.13%"
Set-TextValue "G2" "18"

Set-TextValue "D3" "49.15"
Set-TextValue "E3" "11.16%"
Set-TextValue "G3" "18"

Set-TextValue "D4" "5.310"
Set-TextValue "E4" "4.41%"
Set-TextValue "G4" "18"

Set-TextValue "D5" "0.08068"
Set-TextValue "E5" "4.76%"
Set-TextValue "G5" "18"

Set-TextValue "D6" "4.618"
Set-TextValue "E6" "4.31%"
Set-TextValue "G6" "18"

Set-TextValue "D7" "1.335"
Set-TextValue "E7" "27.29%"
Set-TextValue "G7" "18"

Set-TextValue "E8" "1.96%"
Set-TextValue "G8" "18"

Set-TextValue "D9" "0.1263"
Set-TextValue "E9" "-1.41%"
Set-TextValue "G9" "18"

Set-TextValue "D10" "0.1972"
Set-TextValue "E10" "5.53%"
Set-TextValue "G10" "18"

Set-TextValue "D11" "0.09602"
Set-TextValue "E11" "3.63%"
Set-TextValue "G11" "18"

Set-TextValue "D12" "0.04721"
Set-TextValue "E12" "13.75%"
Set-TextValue "G12" "18"

Set-TextValue "E13" "0.04%"
Set-TextValue "G13" "18"

Set-TextValue "D14" "0.001320"
Set-TextValue "E14" "3.11%"
Set-TextValue "G14" "18"

Set-TextValue "D15" "0.04206"
Set-TextValue "E15" "0.42%"
Set-TextValue "G15" "18"

Set-TextValue "D16" "0.005794"
Set-TextValue "E16" "0.68%"
Set-TextValue "G16" "18"

Set-TextValue "E17" "-0.06%"
Set-TextValue "G17" "18"

Set-TextValue "D18" "2.443"
Set-TextValue "E18" "4.82%"
Set-TextValue "G18" "18"

Set-TextValue "D19" "0.3526"
Set-TextValue "E19" "5.57%"
Set-TextValue "G19" "18"

Set-TextValue "D20" "8.056"
Set-TextValue "E20" "-0.52%"
Set-TextValue "G20" "18"

Set-TextValue "D21" "0.1367"
Set-TextValue "E21" "-2.34%"
Set-TextValue "G21" "18"

Set-TextValue "D22" "0.3077"
Set-TextValue "E22" "-3.18%"
Set-TextValue "G22" "18"

Set-TextValue "D23" "0.001307"
Set-TextValue "E23" "1.74%"
Set-TextValue "G23" "18"

Set-TextValue "D24" "0.004311"
Set-TextValue "E24" "-2.26%"
Set-TextValue "G24" "18"

Set-TextValue "E25" "-0.04%"
Set-TextValue "G25" "18"

Set-TextValue "E26" "-95.27%"
Set-TextValue "G26" "18"

Set-TextValue "G27" "18"

Set-TextValue "G28" "18"

Set-TextValue "G29" "18"

Set-TextValue "G30" "18"

Set-TextValue "G31" "18"

Set-TextValue "G32" "18"

Set-TextValue "G33" "18"

Set-TextValue "G34" "18"

Set-TextValue "G35" "18"

Set-TextValue "G36" "18"

Set-TextValue "G37" "18"

Set-TextValue "D38" "0.02716"
Set-TextValue "E38" "8.60%"
Set-TextValue "G38" "18"

Set-TextValue "D39" "0.06014"
Set-TextValue "E39" "13.40%"
Set-TextValue "G39" "18"

Set-TextValue "E40" "82.93%"
Set-TextValue "G40" "18"

Set-TextValue "D41" "0.008020"
Set-TextValue "E41" "3.75%"
Set-TextValue "G41" "18"

Set-TextValue "E42" "8.69%"
Set-TextValue "G42" "18"

Set-TextValue "D43" "0.007912"
Set-TextValue "E43" "7.64%"
Set-TextValue "G43" "18"

Set-TextValue "D44" "0.007880"
Set-TextValue "E44" "4.79%"
Set-TextValue "G44" "18"

Set-TextValue "D45" "0.3498"
Set-TextValue "E45" "15.71%"
Set-TextValue "G45" "18"

Set-TextValue "D46" "0.00006889"
Set-TextValue "E46" "2.88%"
Set-TextValue "G46" "18"

Set-TextValue "E47" "-0.01%"
Set-TextValue "G47" "18"

Set-TextValue "D48" "0.05951"
Set-TextValue "E48" "37.19%"
Set-TextValue "G48" "18"

Set-TextValue "D49" "0.004000"
Set-TextValue "E49" "-4.77%"
Set-TextValue "G49" "18"

Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "-0.01%"
Set-TextValue "G50" "18"

Set-TextValue "E51" "-0.01%"
Set-TextValue "G51" "18"
